$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 16:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1595318
$ws.Range("C4").Value = 2595
$ws.Range("D4").Value = 370973
$ws.Range("E4").Value = 1129324
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = 95021

# Row 11 - Alemania
$ws.Range("B11").Value = 178671
$ws.Range("C11").Value = 140
$ws.Range("E11").Value = 12400

# Rows 37/38 - Kuwait and Sudafrica swap rank (Sudafrica now above Kuwait)
# and both get refreshed case data.
$ws.Range("A37").Value = "Sudafrica"
$ws.Range("B37").Value = 19137
$ws.Range("C37").Value = 1134
$ws.Range("D37").Value = 8950
$ws.Range("E37").Value = 9818
$ws.Range("G37").Value = 30
$ws.Range("H37").Value = 369

$ws.Range("A38").Value = "Kuwait"
$ws.Range("B38").Value = 18609
$ws.Range("C38").Value = 1041
$ws.Range("D38").Value = 5205
$ws.Range("E38").Value = 13275
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 129
